# Natmi following Dr Hou advice
# Updates the LR-pairs Vegfa-Gpc1 sheet: recomputed communication-score
# columns for the existing ECs/FAPs/M2 sending-cluster rows, and adds the
# previously-missing sCs sending-cluster rows (14-17) completing the
# 4x4 sending x target cluster matrix.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Vegfa"
$ws.Cells.Item(2,3).Value = "Gpc1"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 4.979788333333333
$ws.Cells.Item(2,8).Value = 14.939365
$ws.Cells.Item(2,9).Value = 0.129176854764059
$ws.Cells.Item(2,10).Value = 0.129176854764059
$ws.Cells.Item(2,11).Value = 2
$ws.Cells.Item(2,12).Value = 0.6666666666666666
$ws.Cells.Item(2,13).Value = 1.628421
$ws.Cells.Item(2,14).Value = 4.885263
$ws.Cells.Item(2,15).Value = 0.048329411442081
$ws.Cells.Item(2,16).Value = 0.048329411442081
$ws.Cells.Item(2,17).Value = 8.109191897555
$ws.Cells.Item(2,18).Value = 72.98272707799501
$ws.Cells.Item(2,19).Value = 0.006243041362686147
$ws.Cells.Item(2,20).Value = 0.006243041362686148

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Vegfa"
$ws.Cells.Item(3,3).Value = "Gpc1"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 4.979788333333333
$ws.Cells.Item(3,8).Value = 14.939365
$ws.Cells.Item(3,9).Value = 0.129176854764059
$ws.Cells.Item(3,10).Value = 0.129176854764059
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 9.459065000000001
$ws.Cells.Item(3,14).Value = 28.377195
$ws.Cells.Item(3,15).Value = 0.2807327123897247
$ws.Cells.Item(3,16).Value = 0.2807327123897247
$ws.Cells.Item(3,17).Value = 47.10414153124167
$ws.Cells.Item(3,18).Value = 423.937273781175
$ws.Cells.Item(3,19).Value = 0.03626416881588781
$ws.Cells.Item(3,20).Value = 0.03626416881588781

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Vegfa"
$ws.Cells.Item(4,3).Value = "Gpc1"
$ws.Cells.Item(4,4).Value = "M2"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 4.979788333333333
$ws.Cells.Item(4,8).Value = 14.939365
$ws.Cells.Item(4,9).Value = 0.129176854764059
$ws.Cells.Item(4,10).Value = 0.129176854764059
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 0.6418243333333334
$ws.Cells.Item(4,14).Value = 1.925473
$ws.Cells.Item(4,15).Value = 0.01904850912583786
$ws.Cells.Item(4,16).Value = 0.01904850912583786
$ws.Cells.Item(4,17).Value = 3.196149327182778
$ws.Cells.Item(4,18).Value = 28.765343944645
$ws.Cells.Item(4,19).Value = 0.002460626496820209
$ws.Cells.Item(4,20).Value = 0.00246062649682021

# Row 5
$ws.Cells.Item(5,1).Value = "ECs"
$ws.Cells.Item(5,2).Value = "Vegfa"
$ws.Cells.Item(5,3).Value = "Gpc1"
$ws.Cells.Item(5,4).Value = "sCs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 4.979788333333333
$ws.Cells.Item(5,8).Value = 14.939365
$ws.Cells.Item(5,9).Value = 0.129176854764059
$ws.Cells.Item(5,10).Value = 0.129176854764059
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 21.96489266666667
$ws.Cells.Item(5,14).Value = 65.894678
$ws.Cells.Item(5,15).Value = 0.6518893670423564
$ws.Cells.Item(5,16).Value = 0.6518893670423563
$ws.Cells.Item(5,17).Value = 109.3805162443856
$ws.Cells.Item(5,18).Value = 984.42464619947
$ws.Cells.Item(5,19).Value = 0.08420901808866481
$ws.Cells.Item(5,20).Value = 0.08420901808866481

# Row 6
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Vegfa"
$ws.Cells.Item(6,3).Value = "Gpc1"
$ws.Cells.Item(6,4).Value = "ECs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 18.019504
$ws.Cells.Item(6,8).Value = 54.058512
$ws.Cells.Item(6,9).Value = 0.467430078412646
$ws.Cells.Item(6,10).Value = 0.4674300784126461
$ws.Cells.Item(6,11).Value = 2
$ws.Cells.Item(6,12).Value = 0.6666666666666666
$ws.Cells.Item(6,13).Value = 1.628421
$ws.Cells.Item(6,14).Value = 4.885263
$ws.Cells.Item(6,15).Value = 0.048329411442081
$ws.Cells.Item(6,16).Value = 0.048329411442081
$ws.Cells.Item(6,17).Value = 29.343338723184
$ws.Cells.Item(6,18).Value = 264.090048508656
$ws.Cells.Item(6,19).Value = 0.02259062058000895
$ws.Cells.Item(6,20).Value = 0.02259062058000896

# Row 7
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Vegfa"
$ws.Cells.Item(7,3).Value = "Gpc1"
$ws.Cells.Item(7,4).Value = "FAPs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 18.019504
$ws.Cells.Item(7,8).Value = 54.058512
$ws.Cells.Item(7,9).Value = 0.467430078412646
$ws.Cells.Item(7,10).Value = 0.4674300784126461
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 9.459065000000001
$ws.Cells.Item(7,14).Value = 28.377195
$ws.Cells.Item(7,15).Value = 0.2807327123897247
$ws.Cells.Item(7,16).Value = 0.2807327123897247
$ws.Cells.Item(7,17).Value = 170.44765960376
$ws.Cells.Item(7,18).Value = 1534.02893643384
$ws.Cells.Item(7,19).Value = 0.1312229137653238
$ws.Cells.Item(7,20).Value = 0.1312229137653238

# Row 8
$ws.Cells.Item(8,1).Value = "FAPs"
$ws.Cells.Item(8,2).Value = "Vegfa"
$ws.Cells.Item(8,3).Value = "Gpc1"
$ws.Cells.Item(8,4).Value = "M2"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 18.019504
$ws.Cells.Item(8,8).Value = 54.058512
$ws.Cells.Item(8,9).Value = 0.467430078412646
$ws.Cells.Item(8,10).Value = 0.4674300784126461
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 0.6418243333333334
$ws.Cells.Item(8,14).Value = 1.925473
$ws.Cells.Item(8,15).Value = 0.01904850912583786
$ws.Cells.Item(8,16).Value = 0.01904850912583786
$ws.Cells.Item(8,17).Value = 11.56535614179734
$ws.Cells.Item(8,18).Value = 104.088205276176
$ws.Cells.Item(8,19).Value = 0.008903846114334393
$ws.Cells.Item(8,20).Value = 0.008903846114334394

# Row 9
$ws.Cells.Item(9,1).Value = "FAPs"
$ws.Cells.Item(9,2).Value = "Vegfa"
$ws.Cells.Item(9,3).Value = "Gpc1"
$ws.Cells.Item(9,4).Value = "sCs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 18.019504
$ws.Cells.Item(9,8).Value = 54.058512
$ws.Cells.Item(9,9).Value = 0.467430078412646
$ws.Cells.Item(9,10).Value = 0.4674300784126461
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 21.96489266666667
$ws.Cells.Item(9,14).Value = 65.894678
$ws.Cells.Item(9,15).Value = 0.6518893670423564
$ws.Cells.Item(9,16).Value = 0.6518893670423563
$ws.Cells.Item(9,17).Value = 395.7964712665707
$ws.Cells.Item(9,18).Value = 3562.168241399136
$ws.Cells.Item(9,19).Value = 0.3047126979529788
$ws.Cells.Item(9,20).Value = 0.3047126979529788

# Row 10
$ws.Cells.Item(10,1).Value = "M2"
$ws.Cells.Item(10,2).Value = "Vegfa"
$ws.Cells.Item(10,3).Value = "Gpc1"
$ws.Cells.Item(10,4).Value = "ECs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 8.752692000000001
$ws.Cells.Item(10,8).Value = 26.258076
$ws.Cells.Item(10,9).Value = 0.2270468436801446
$ws.Cells.Item(10,10).Value = 0.2270468436801446
$ws.Cells.Item(10,11).Value = 2
$ws.Cells.Item(10,12).Value = 0.6666666666666666
$ws.Cells.Item(10,13).Value = 1.628421
$ws.Cells.Item(10,14).Value = 4.885263
$ws.Cells.Item(10,15).Value = 0.048329411442081
$ws.Cells.Item(10,16).Value = 0.048329411442081
$ws.Cells.Item(10,17).Value = 14.253067459332
$ws.Cells.Item(10,18).Value = 128.277607133988
$ws.Cells.Item(10,19).Value = 0.01097304032484355
$ws.Cells.Item(10,20).Value = 0.01097304032484355

# Row 11
$ws.Cells.Item(11,1).Value = "M2"
$ws.Cells.Item(11,2).Value = "Vegfa"
$ws.Cells.Item(11,3).Value = "Gpc1"
$ws.Cells.Item(11,4).Value = "FAPs"
$ws.Cells.Item(11,5).Value = 3
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 8.752692000000001
$ws.Cells.Item(11,8).Value = 26.258076
$ws.Cells.Item(11,9).Value = 0.2270468436801446
$ws.Cells.Item(11,10).Value = 0.2270468436801446
$ws.Cells.Item(11,11).Value = 3
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 9.459065000000001
$ws.Cells.Item(11,14).Value = 28.377195
$ws.Cells.Item(11,15).Value = 0.2807327123897247
$ws.Cells.Item(11,16).Value = 0.2807327123897247
$ws.Cells.Item(11,17).Value = 82.79228255298003
$ws.Cells.Item(11,18).Value = 745.13054297682
$ws.Cells.Item(11,19).Value = 0.06373947626585282
$ws.Cells.Item(11,20).Value = 0.06373947626585282

# Row 12
$ws.Cells.Item(12,1).Value = "M2"
$ws.Cells.Item(12,2).Value = "Vegfa"
$ws.Cells.Item(12,3).Value = "Gpc1"
$ws.Cells.Item(12,4).Value = "M2"
$ws.Cells.Item(12,5).Value = 3
$ws.Cells.Item(12,6).Value = 1
$ws.Cells.Item(12,7).Value = 8.752692000000001
$ws.Cells.Item(12,8).Value = 26.258076
$ws.Cells.Item(12,9).Value = 0.2270468436801446
$ws.Cells.Item(12,10).Value = 0.2270468436801446
$ws.Cells.Item(12,11).Value = 3
$ws.Cells.Item(12,12).Value = 1
$ws.Cells.Item(12,13).Value = 0.6418243333333334
$ws.Cells.Item(12,14).Value = 1.925473
$ws.Cells.Item(12,15).Value = 0.01904850912583786
$ws.Cells.Item(12,16).Value = 0.01904850912583786
$ws.Cells.Item(12,17).Value = 5.617690707772002
$ws.Cells.Item(12,18).Value = 50.55921636994801
$ws.Cells.Item(12,19).Value = 0.004324903873833916
$ws.Cells.Item(12,20).Value = 0.004324903873833916

# Row 13
$ws.Cells.Item(13,1).Value = "M2"
$ws.Cells.Item(13,2).Value = "Vegfa"
$ws.Cells.Item(13,3).Value = "Gpc1"
$ws.Cells.Item(13,4).Value = "sCs"
$ws.Cells.Item(13,5).Value = 3
$ws.Cells.Item(13,6).Value = 1
$ws.Cells.Item(13,7).Value = 8.752692000000001
$ws.Cells.Item(13,8).Value = 26.258076
$ws.Cells.Item(13,9).Value = 0.2270468436801446
$ws.Cells.Item(13,10).Value = 0.2270468436801446
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = 21.96489266666667
$ws.Cells.Item(13,14).Value = 65.894678
$ws.Cells.Item(13,15).Value = 0.6518893670423564
$ws.Cells.Item(13,16).Value = 0.6518893670423563
$ws.Cells.Item(13,17).Value = 192.2519403243921
$ws.Cells.Item(13,18).Value = 1730.267462919528
$ws.Cells.Item(13,19).Value = 0.1480094232156143
$ws.Cells.Item(13,20).Value = 0.1480094232156143

# Row 14
$ws.Cells.Item(14,1).Value = "sCs"
$ws.Cells.Item(14,2).Value = "Vegfa"
$ws.Cells.Item(14,3).Value = "Gpc1"
$ws.Cells.Item(14,4).Value = "ECs"
$ws.Cells.Item(14,5).Value = 3
$ws.Cells.Item(14,6).Value = 1
$ws.Cells.Item(14,7).Value = 6.798175000000001
$ws.Cells.Item(14,8).Value = 20.394525
$ws.Cells.Item(14,9).Value = 0.1763462231431503
$ws.Cells.Item(14,10).Value = 0.1763462231431503
$ws.Cells.Item(14,11).Value = 2
$ws.Cells.Item(14,12).Value = 0.6666666666666666
$ws.Cells.Item(14,13).Value = 1.628421
$ws.Cells.Item(14,14).Value = 4.885263
$ws.Cells.Item(14,15).Value = 0.048329411442081
$ws.Cells.Item(14,16).Value = 0.048329411442081
$ws.Cells.Item(14,17).Value = 11.070290931675
$ws.Cells.Item(14,18).Value = 99.63261838507501
$ws.Cells.Item(14,19).Value = 0.008522709174542339
$ws.Cells.Item(14,20).Value = 0.00852270917454234

# Row 15
$ws.Cells.Item(15,1).Value = "sCs"
$ws.Cells.Item(15,2).Value = "Vegfa"
$ws.Cells.Item(15,3).Value = "Gpc1"
$ws.Cells.Item(15,4).Value = "FAPs"
$ws.Cells.Item(15,5).Value = 3
$ws.Cells.Item(15,6).Value = 1
$ws.Cells.Item(15,7).Value = 6.798175000000001
$ws.Cells.Item(15,8).Value = 20.394525
$ws.Cells.Item(15,9).Value = 0.1763462231431503
$ws.Cells.Item(15,10).Value = 0.1763462231431503
$ws.Cells.Item(15,11).Value = 3
$ws.Cells.Item(15,12).Value = 1
$ws.Cells.Item(15,13).Value = 9.459065000000001
$ws.Cells.Item(15,14).Value = 28.377195
$ws.Cells.Item(15,15).Value = 0.2807327123897247
$ws.Cells.Item(15,16).Value = 0.2807327123897247
$ws.Cells.Item(15,17).Value = 64.304379206375
$ws.Cells.Item(15,18).Value = 578.7394128573751
$ws.Cells.Item(15,19).Value = 0.04950615354266024
$ws.Cells.Item(15,20).Value = 0.04950615354266024

# Row 16
$ws.Cells.Item(16,1).Value = "sCs"
$ws.Cells.Item(16,2).Value = "Vegfa"
$ws.Cells.Item(16,3).Value = "Gpc1"
$ws.Cells.Item(16,4).Value = "M2"
$ws.Cells.Item(16,5).Value = 3
$ws.Cells.Item(16,6).Value = 1
$ws.Cells.Item(16,7).Value = 6.798175000000001
$ws.Cells.Item(16,8).Value = 20.394525
$ws.Cells.Item(16,9).Value = 0.1763462231431503
$ws.Cells.Item(16,10).Value = 0.1763462231431503
$ws.Cells.Item(16,11).Value = 3
$ws.Cells.Item(16,12).Value = 1
$ws.Cells.Item(16,13).Value = 0.6418243333333334
$ws.Cells.Item(16,14).Value = 1.925473
$ws.Cells.Item(16,15).Value = 0.01904850912583786
$ws.Cells.Item(16,16).Value = 0.01904850912583786
$ws.Cells.Item(16,17).Value = 4.363234137258335
$ws.Cells.Item(16,18).Value = 39.26910723532501
$ws.Cells.Item(16,19).Value = 0.003359132640849338
$ws.Cells.Item(16,20).Value = 0.003359132640849339

# Row 17
$ws.Cells.Item(17,1).Value = "sCs"
$ws.Cells.Item(17,2).Value = "Vegfa"
$ws.Cells.Item(17,3).Value = "Gpc1"
$ws.Cells.Item(17,4).Value = "sCs"
$ws.Cells.Item(17,5).Value = 3
$ws.Cells.Item(17,6).Value = 1
$ws.Cells.Item(17,7).Value = 6.798175000000001
$ws.Cells.Item(17,8).Value = 20.394525
$ws.Cells.Item(17,9).Value = 0.1763462231431503
$ws.Cells.Item(17,10).Value = 0.1763462231431503
$ws.Cells.Item(17,11).Value = 3
$ws.Cells.Item(17,12).Value = 1
$ws.Cells.Item(17,13).Value = 21.96489266666667
$ws.Cells.Item(17,14).Value = 65.894678
$ws.Cells.Item(17,15).Value = 0.6518893670423564
$ws.Cells.Item(17,16).Value = 0.6518893670423563
$ws.Cells.Item(17,17).Value = 149.3211842042167
$ws.Cells.Item(17,18).Value = 1343.89065783795
$ws.Cells.Item(17,19).Value = 0.1149582277850984
$ws.Cells.Item(17,20).Value = 0.1149582277850984

